$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (9) to hold the new
# "Distribution channel code" field; this pushes the existing
# "Budget" column from I to J.
$ws.Columns.Item(9).Insert()

# Header for the new column.
$ws.Cells.Item(1, 9).Value = "Distribution channel code"

# Distribution channel code values for the data rows.
$ws.Cells.Item(2, 9).Value = "TR"
$ws.Cells.Item(3, 9).Value = "GO"

# Match the new column's width from the diff (22.5546875 characters).
$ws.Columns.Item(9).ColumnWidth = 21.6666667

# Update the active selection to reflect the shifted columns.
$ws.Range("M13").Select()
